$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("S2LAR_Sedan_HambaLG_f")
Write-Host $ws1.Name
